$d = $word.ActiveDocument
$t = $d.Tables(1)

# 50÷9= -> 79÷7=
$cell = $t.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "50÷9=") {
    throw "Unexpected cell text at Row 1 Col 1: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "79÷7="

# 49÷2= -> 19÷6=
$cell = $t.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "49÷2=") {
    throw "Unexpected cell text at Row 1 Col 2: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "19÷6="

# 90÷9= -> 10÷5=
$cell = $t.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "90÷9=") {
    throw "Unexpected cell text at Row 1 Col 3: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "10÷5="

# 49÷7= -> 46÷4=
$cell = $t.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "49÷7=") {
    throw "Unexpected cell text at Row 1 Col 4: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "46÷4="

# 23÷7= -> 97÷5=
$cell = $t.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "23÷7=") {
    throw "Unexpected cell text at Row 1 Col 5: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "97÷5="

# 73÷9= -> 96÷6=
$cell = $t.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "73÷9=") {
    throw "Unexpected cell text at Row 5 Col 1: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "96÷6="

# 52÷4= -> 39÷7=
$cell = $t.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "52÷4=") {
    throw "Unexpected cell text at Row 5 Col 2: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "39÷7="

# 47÷5= -> 89÷7=
$cell = $t.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "47÷5=") {
    throw "Unexpected cell text at Row 5 Col 3: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "89÷7="

# 71÷5= -> 97÷2=
$cell = $t.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "71÷5=") {
    throw "Unexpected cell text at Row 5 Col 4: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "97÷2="

# 60÷8= -> 41÷4=
$cell = $t.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "60÷8=") {
    throw "Unexpected cell text at Row 5 Col 5: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "41÷4="

# 11÷2= -> 91÷8=
$cell = $t.Cell(9, 1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "11÷2=") {
    throw "Unexpected cell text at Row 9 Col 1: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "91÷8="

# 84÷7= -> 81÷5=
$cell = $t.Cell(9, 2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "84÷7=") {
    throw "Unexpected cell text at Row 9 Col 2: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "81÷5="

# 46÷9= -> 21÷6=
$cell = $t.Cell(9, 3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "46÷9=") {
    throw "Unexpected cell text at Row 9 Col 3: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "21÷6="

# 31÷4= -> 18÷3=
$cell = $t.Cell(9, 4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "31÷4=") {
    throw "Unexpected cell text at Row 9 Col 4: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "18÷3="

# 78÷9= -> 78÷6=
$cell = $t.Cell(9, 5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "78÷9=") {
    throw "Unexpected cell text at Row 9 Col 5: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "78÷6="

# 74÷3= -> 38÷2=
$cell = $t.Cell(13, 1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "74÷3=") {
    throw "Unexpected cell text at Row 13 Col 1: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "38÷2="

# 35÷3= -> 48÷9=
$cell = $t.Cell(13, 2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "35÷3=") {
    throw "Unexpected cell text at Row 13 Col 2: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "48÷9="

# 40÷6= -> 57÷9=
$cell = $t.Cell(13, 3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "40÷6=") {
    throw "Unexpected cell text at Row 13 Col 3: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "57÷9="

# 10÷6= -> 72÷8=
$cell = $t.Cell(13, 4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "10÷6=") {
    throw "Unexpected cell text at Row 13 Col 4: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "72÷8="

# 42÷5= -> 90÷8=
$cell = $t.Cell(13, 5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "42÷5=") {
    throw "Unexpected cell text at Row 13 Col 5: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "90÷8="

# 17÷8= -> 78÷8=
$cell = $t.Cell(17, 1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "17÷8=") {
    throw "Unexpected cell text at Row 17 Col 1: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "78÷8="

# 58÷9= -> 60÷5=
$cell = $t.Cell(17, 2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "58÷9=") {
    throw "Unexpected cell text at Row 17 Col 2: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "60÷5="

# 20÷6= -> 11÷6=
$cell = $t.Cell(17, 3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "20÷6=") {
    throw "Unexpected cell text at Row 17 Col 3: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "11÷6="

# 81÷5= -> 92÷9=
$cell = $t.Cell(17, 4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "81÷5=") {
    throw "Unexpected cell text at Row 17 Col 4: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "92÷9="

# 74÷6= -> 20÷8=
$cell = $t.Cell(17, 5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "74÷6=") {
    throw "Unexpected cell text at Row 17 Col 5: [" + $cell.Range.Text + "]"
}
$cell.Range.Text = "20÷8="
